$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:F2").NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = "87"
$ws.Cells.Item(2, 4).Value = "61"
$ws.Cells.Item(2, 5).Value = "10"
$ws.Cells.Item(2, 6).Value = "4"

$ws.Range("C4:F4").NumberFormat = "@"
$ws.Cells.Item(4, 3).Value = "5"
$ws.Cells.Item(4, 4).Value = "6"
$ws.Cells.Item(4, 5).Value = "1"
$ws.Cells.Item(4, 6).Value = "0"

$ws.Range("C5:E5").NumberFormat = "@"
$ws.Cells.Item(5, 3).Value = "9"
$ws.Cells.Item(5, 4).Value = "10"
$ws.Cells.Item(5, 5).Value = "1"

$ws.Range("C6:F6").NumberFormat = "@"
$ws.Cells.Item(6, 3).Value = "0"
$ws.Cells.Item(6, 4).Value = "1"
$ws.Cells.Item(6, 5).Value = "0"
$ws.Cells.Item(6, 6).Value = "0"

$ws.Range("C7:F7").NumberFormat = "@"
$ws.Cells.Item(7, 3).Value = "2"
$ws.Cells.Item(7, 4).Value = "4"
$ws.Cells.Item(7, 5).Value = "0"
$ws.Cells.Item(7, 6).Value = "0"

$ws.Range("C8:E8").NumberFormat = "@"
$ws.Cells.Item(8, 3).Value = "81"
$ws.Cells.Item(8, 4).Value = "53"
$ws.Cells.Item(8, 5).Value = "13"

$ws.Range("C9:F9").NumberFormat = "@"
$ws.Cells.Item(9, 3).Value = "9"
$ws.Cells.Item(9, 4).Value = "14"
$ws.Cells.Item(9, 5).Value = "1"
$ws.Cells.Item(9, 6).Value = "0"

$ws.Range("C10:F10").NumberFormat = "@"
$ws.Cells.Item(10, 3).Value = "29"
$ws.Cells.Item(10, 4).Value = "20"
$ws.Cells.Item(10, 5).Value = "3"
$ws.Cells.Item(10, 6).Value = "1"

$ws.Range("C11:E11").NumberFormat = "@"
$ws.Cells.Item(11, 3).Value = "26"
$ws.Cells.Item(11, 4).Value = "13"
$ws.Cells.Item(11, 5).Value = "6"

$ws.Range("C12:F12").NumberFormat = "@"
$ws.Cells.Item(12, 3).Value = "24"
$ws.Cells.Item(12, 4).Value = "18"
$ws.Cells.Item(12, 5).Value = "2"
$ws.Cells.Item(12, 6).Value = "1"

$ws.Range("C13:F13").NumberFormat = "@"
$ws.Cells.Item(13, 3).Value = "0"
$ws.Cells.Item(13, 4).Value = "1"
$ws.Cells.Item(13, 5).Value = "0"
$ws.Cells.Item(13, 6).Value = "0"

$ws.Range("C14:F14").NumberFormat = "@"
$ws.Cells.Item(14, 3).Value = "58"
$ws.Cells.Item(14, 4).Value = "35"
$ws.Cells.Item(14, 5).Value = "4"
$ws.Cells.Item(14, 6).Value = "4"

$ws.Range("C15:F15").NumberFormat = "@"
$ws.Cells.Item(15, 3).Value = "22"
$ws.Cells.Item(15, 4).Value = "17"
$ws.Cells.Item(15, 5).Value = "2"
$ws.Cells.Item(15, 6).Value = "1"

